$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update headers
$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Data rows (Índice, Distancia(avg), max, min, Tempo)
$data = @(
    @(0, 1283.133333333333, 1440, 1223, 0.02861672242482503),
    @(1, 1327.433333333333, 1436, 1272, 0.02868072191874186),
    @(2, 1304.6,            1451, 1224, 0.03137623469034831),
    @(3, 1313.033333333333, 1438, 1272, 0.02846802075703939),
    @(4, 1285.366666666667, 1389, 1221, 0.0287573496500651),
    @(5, 1178,              1316, 1068, 0.02873221238454183),
    @(6, 1147.9,            1365, 1057, 0.02886221408843994),
    @(7, 1131.466666666667, 1220, 1068, 0.02888691425323486),
    @(8, 1260.9,            1375, 1221, 0.02911811669667562),
    @(9, 1152.9,            1390, 1057, 0.02915050188700358)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $row++
}
